$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 updates ---------------------------------------------------------
# B9: TrainingAndTestDataScenario1 -> TrainingAndTestDataScenario3
$ws.Range("B9").Value = "TrainingAndTestDataScenario3"

# D9: (blank) -> 180  (E9 formula =D9/60/60 recalculates automatically)
$ws.Range("D9").Value = 180

$ws.Range("F9").Value = "^"
$ws.Range("G9").Value = "^"
$ws.Range("H9").Value = "^"
$ws.Range("J9").Value = "^"
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = "^"
$ws.Range("M9").Value = 0.313
$ws.Range("N9").Value = 0.255
$ws.Range("O9").Value = "7/30/2023"

# --- Row 10: new data row ---------------------------------------------------
# New shared strings are appended in the order they are first assigned, so
# write the new-string cells in the same order the source sheet introduces
# them: A10, C10, I10, then I9.
$ws.Range("A10").Value = "HyperparameterSweepResultsScenario8"
$ws.Range("B10").Value = "^"
$ws.Range("C10").Value = "NeuralNetworkScenario5"
$ws.Range("D10").Value = 7200
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = "^"
$ws.Range("H10").Value = "^"
$ws.Range("I10").Value = "[0.001,0.01]"
$ws.Range("J10").Value = "^"
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = "^"
$ws.Range("M10").Value = 0.78749999999999998
$ws.Range("N10").Value = 0.69499999999999995
$ws.Range("O10").Value = "7/31/2023"

$ws.Range("I9").Value = "[0.00001,0.0005]"

# --- Selection follows the edits, ending on M11 -----------------------------
$ws.Range("M11").Select() | Out-Null
